$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05678866666666666
$ws.Range("H2").Value = 0.170366
$ws.Range("I2").Value = 0.1297730042656916
$ws.Range("J2").Value = 0.1297730042656917
$ws.Range("M2").Value = 0.4010836666666667
$ws.Range("O2").Value = 0.0699907978100882
$ws.Range("P2").Value = 0.0699907978100882
$ws.Range("Q2").Value = 0.02277700665177778
$ws.Range("R2").Value = 0.204993059866
$ws.Range("S2").Value = 0.009082916102767737
$ws.Range("T2").Value = 0.009082916102767739

$ws.Range("G3").Value = 0.05678866666666666
$ws.Range("H3").Value = 0.170366
$ws.Range("I3").Value = 0.1297730042656916
$ws.Range("J3").Value = 0.1297730042656917
$ws.Range("M3").Value = 5.329436333333334
$ws.Range("N3").Value = 15.988309
$ws.Range("O3").Value = 0.9300092021899118
$ws.Range("P3").Value = 0.9300092021899117
$ws.Range("Q3").Value = 0.3026515834548889
$ws.Range("R3").Value = 2.723864251094
$ws.Range("S3").Value = 0.1206900881629239
$ws.Range("T3").Value = 0.1206900881629239

$ws.Range("I4").Value = 0.8702269957343084
$ws.Range("J4").Value = 0.8702269957343083
$ws.Range("M4").Value = 0.4010836666666667
$ws.Range("O4").Value = 0.0699907978100882
$ws.Range("P4").Value = 0.0699907978100882
$ws.Range("Q4").Value = 0.1527372058815556
$ws.Range("S4").Value = 0.06090788170732046
$ws.Range("T4").Value = 0.06090788170732046

$ws.Range("I5").Value = 0.8702269957343084
$ws.Range("J5").Value = 0.8702269957343083
$ws.Range("M5").Value = 5.329436333333334
$ws.Range("N5").Value = 15.988309
$ws.Range("O5").Value = 0.9300092021899118
$ws.Range("P5").Value = 0.9300092021899117
$ws.Range("Q5").Value = 2.029509756011778
$ws.Range("R5").Value = 18.265587804106
$ws.Range("S5").Value = 0.809319114026988
$ws.Range("T5").Value = 0.8093191140269878
